$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D1:R1 hold values that look numeric/date-like to Excel's parser
# (dates, percentages, zero-padded decimals, etc.) but must be preserved
# as literal text, matching the source data. Force text format before
# assigning so Excel does not auto-convert them.
$ws.Range("D1:R1").NumberFormat = "@"

$ws.Range("A1").Value = "Swoopes"
$ws.Range("B1").Value = "Tyrone"
$ws.Range("C1").Value = "TE"
$ws.Range("D1").Value = "2018-10-14"
$ws.Range("E1").Value = "6"
$ws.Range("F1").Value = "23.334"
$ws.Range("G1").Value = "SEA"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "OAK"
$ws.Range("J1").Value = "W 27-3"
$ws.Range("K1").Value = "*"
$ws.Range("L1").Value = "1"
$ws.Range("M1").Value = "1"
$ws.Range("N1").Value = "23"
$ws.Range("O1").Value = "23.00"
$ws.Range("P1").Value = "0"
$ws.Range("Q1").Value = "100.0%"
$ws.Range("R1").Value = "23.00"
$ws.Range("S1").Value = 2.3
